# Apply "break out stock.yaml completed" edit:
#  1. Fix D245:D248 on the "day" sheet - these were written as text
#     ("inlineStr") by mistake; they should be plain numbers (BSE codes).
#  2. Append 14 new rows (117-130) of weekly stock data to the "week" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fix bsecode column (D) for rows 245-248 on the "day" sheet so the
#    values are stored as numbers instead of text.
# ---------------------------------------------------------------------
$dayWs = $wb.Worksheets.Item("day")

$dayWs.Range("D245").Value = 540115
$dayWs.Range("D246").Value = 500820
$dayWs.Range("D247").Value = 533398
$dayWs.Range("D248").Value = 532187

# ---------------------------------------------------------------------
# 2) Append new rows to the "week" sheet.
# ---------------------------------------------------------------------
$weekWs = $wb.Worksheets.Item("week")

$newRows = @(
    @(1,  "BOSCHLTD",    "Bosch Limited",                                   500530, -1.34, 33909.75,   28194,     "week", "02/08/2024 11:33:12"),
    @(2,  "POLYCAB",     "Polycab India Ltd",                               542652, -0.35, 6751.5,     285579,    "week", "02/08/2024 11:33:12"),
    @(3,  "JKCEMENT",    "Jk Cement Limited",                               532644, -3.1,  4284.6,     98277,     "week", "02/08/2024 11:33:12"),
    @(4,  "CUMMINSIND",  "Cummins India Limited",                           500480, -7.97, 3510.85,    4800902,   "week", "02/08/2024 11:33:12"),
    @(5,  "SHRIRAMFIN",  "Shriram Finance Ltd",                             511218, 0.1,   2992.1,     1781076,   "week", "02/08/2024 11:33:12"),
    @(6,  "BHARTIARTL",  "Bharti Airtel Limited",                           532454, -0.78, 1493.8,     3335916,   "week", "02/08/2024 11:33:12"),
    @(7,  "IPCALAB",     "Ipca Laboratories Limited",                       524494, -0.34, 1300.3,     714748,    "week", "02/08/2024 11:33:12"),
    @(8,  "RAMCOCEM",    "The Ramco Cements Limited",                       500260, 0.32,  827.75,     1755087,   "week", "02/08/2024 11:33:12"),
    @(9,  "AUBANK",      "AU Small Finance Bank",                           540611, -1.11, 637.6,      3036465,   "week", "02/08/2024 11:33:12"),
    @(10, "INDHOTEL",    "The Indian Hotels Company Limited",               500850, -0.53, 628.8,      1814826,   "week", "02/08/2024 11:33:12"),
    @(11, "EXIDEIND",    "Exide Industries Limited",                        500086, -0.76, 508.85,     3858656,   "week", "02/08/2024 11:33:12"),
    @(12, "M&MFIN",      "Mahindra & Mahindra Financial Services Limited",  532720, -0.33, 305.45,     1532304,   "week", "02/08/2024 11:33:12"),
    @(13, "RBLBANK",     "Rbl Bank Limited",                                540065, -1.83, 226.9,      5849428,   "week", "02/08/2024 11:33:12"),
    @(14, "IDEA",        "Idea Cellular Limited",                           532822, 0.81,  16.12,      475279605, "week", "02/08/2024 11:33:12")
)

$startRow = 117
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $weekWs.Cells.Item($r, 1).Value = $row[0]
    $weekWs.Cells.Item($r, 2).Value = $row[1]
    $weekWs.Cells.Item($r, 3).Value = $row[2]
    $weekWs.Cells.Item($r, 4).Value = $row[3]
    $weekWs.Cells.Item($r, 5).Value = $row[4]
    $weekWs.Cells.Item($r, 6).Value = $row[5]
    $weekWs.Cells.Item($r, 7).Value = $row[6]
    $weekWs.Cells.Item($r, 8).Value = $row[7]
    $weekWs.Cells.Item($r, 9).Value = $row[8]
}
